$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Work from the bottom of the document upwards so that paragraph
# indices used below stay valid as new paragraphs are inserted.
# ------------------------------------------------------------------

# --- Insert "Other visualizations..." / "Going to switch gears..." ---
# These two new bullet paragraphs land right before the existing
# "Only a single visualization..." paragraph (#25), inheriting its
# ListParagraph / numId=2 / ilvl=0 formatting, and are NOT bold.
$pOnly = $d.Paragraphs(25)
$rOnly = $pOnly.Range
$rOnly.Collapse(1)
$rOnly.InsertParagraphBefore()
$rOnly.InsertParagraphBefore()
$pOther = $d.Paragraphs(25)
$pSwitch = $d.Paragraphs(26)
$pOther.Range.Text = "Other visualizations focused on how hard it is to quantify performance by nation."
$pSwitch.Range.Text = "Going to switch gears with this visualization."

# --- Insert "DATA TIPS" bullet after "For nations competing..." (#21) ---
$pNations = $d.Paragraphs(21)
$rNations = $pNations.Range
$rNations.Collapse(0)
$rNations.InsertParagraphAfter()
$pDataTips = $d.Paragraphs(22)
$pDataTips.Range.Text = "DATA TIPS"
$pDataTips.Range.Font.Bold = 1

# --- Insert "Hover Tools" bullet after "First image provides..." (#7) ---
$pFirstImage = $d.Paragraphs(7)
$rFirstImage = $pFirstImage.Range
$rFirstImage.Collapse(0)
$rFirstImage.InsertParagraphAfter()
$pHover = $d.Paragraphs(8)
$pHover.Range.Text = "Hover Tools"
$pHover.Range.ListFormat.ListLevelNumber = 1
$pHover.Range.Font.Bold = 1

# --- "Visual Framework" (#5): bold existing runs, append new bold run ---
$pVisual = $d.Paragraphs(5)
$pVisual.Range.Font.Bold = 1
$rVisual = $pVisual.Range
$rVisual.Collapse(0)
$rVisual.MoveEnd(1, -1)
$rVisual.InsertAfter(" – Persistent Navigation")
$rVisual.Font.Bold = 1

# --- "Home Page:" (#3) -> "Home Page", centered + bold ---
$d.Content.Find.Execute("Home Page:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Home Page", 2)
$pHome = $d.Paragraphs(3)
$pHome.Alignment = 1
$pHome.Range.Font.Bold = 1

# --- "Presentation Notes" (#1): centered, bold, 14pt (sz/szCs 28) ---
$pTitle = $d.Paragraphs(1)
$pTitle.Alignment = 1
$pTitle.Range.Font.Bold = 1
$pTitle.Range.Font.Size = 14
$pTitle.Range.Font.SizeBi = 14

Write-Output "done"
